$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A3"   = -21.933
    "D18"  = -8.349
    "A21"  = -20.645
    "A23"  = -20.808
    "B24"  = 5.342
    "A25"  = -21.589
    "B28"  = 5.468000000000001
    "B36"  = 7.363
    "B45"  = 5.647
    "B48"  = 5.609
    "B49"  = 6.237
    "D51"  = -8.231
    "B52"  = 5.968000000000001
    "A53"  = -20.335
    "B53"  = 8.151
    "B54"  = 4.886
    "D55"  = -8.349
    "A57"  = -22.291
    "A59"  = -22.401
    "D64"  = -7.407000000000001
    "A69"  = -21.441
    "B70"  = 5.052999999999999
    "A79"  = -21.305
    "D80"  = -8.166
    "A83"  = -21.976
    "B86"  = 5.197
    "B87"  = 4.718
    "D92"  = -7.355000000000001
    "A93"  = -21.533
    "D94"  = -7.703999999999999
    "D96"  = -7.738999999999999
    "B101" = 5.217000000000001
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
